$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text stays literal text, matching the source
# file where every populated data cell is stored as t="str" rather than
# a real number (the "@" text format blocks Excel's automatic numeric
# coercion of values such as "1.643" or "0").
$ws.Range("A9:G24").NumberFormat = "@"
$ws.Range("H9:P24").NumberFormat = "@"
$ws.Range("A25").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("H25:P25").NumberFormat = "@"

# Row 9
$ws.Range("A9").Value = "TQ-17"
$ws.Range("B9").Value = "MERCADO INTERNO"
$ws.Range("C9").Value = "VIBRA PRESIDENTE PRUDENTE"
$ws.Range("D9").Value = "0169/23F"
$ws.Range("E9").Value = "VIBRA CUB TRANSF CONF NF 1277460-0 14.09.23"
$ws.Range("F9").Value = "OLEO DIESEL S10"
$ws.Range("G9").Value = "ONU 1202, DIESEL S10"
$ws.Range("H9").Value = "1.643"
$ws.Range("I9").Value = "2.980"
$ws.Range("J9").Value = "0"
$ws.Range("K9").Value = "0"
$ws.Range("L9").Value = "0"
$ws.Range("M9").Value = "0"
$ws.Range("N9").Value = "1.643"
$ws.Range("O9").Value = "2.980"
$ws.Range("P9").Value = "0"

# Row 10
$ws.Range("A10").Value = "TQ-17"
$ws.Range("B10").Value = "MERCADO INTERNO"
$ws.Range("C10").Value = "VIBRA CUBATAO"
$ws.Range("D10").Value = "0169/23G"
$ws.Range("E10").Value = "VIBRA PRES TRANSF CONF NF 607274-0 26.09.2023"
$ws.Range("F10").Value = "OLEO DIESEL S10"
$ws.Range("G10").Value = "ONU 1202, DIESEL S10"
$ws.Range("H10").Value = "20.556"
$ws.Range("I10").Value = "17.730"
$ws.Range("J10").Value = "0"
$ws.Range("K10").Value = "0"
$ws.Range("L10").Value = "0"
$ws.Range("M10").Value = "0"
$ws.Range("N10").Value = "20.556"
$ws.Range("O10").Value = "17.730"
$ws.Range("P10").Value = "0"

# Row 11
$ws.Range("A11").Value = "TQ-17"
$ws.Range("B11").Value = "MERCADO INTERNO"
$ws.Range("C11").Value = "VIBRA CUBATAO"
$ws.Range("D11").Value = "0177/23B"
$ws.Range("E11").Value = "TRANSF NIMOFAST X VIBRA CUB CONF NF 2217-1 09.10.23"
$ws.Range("F11").Value = "OLEO DIESEL S10"
$ws.Range("G11").Value = "OLEO DIESEL S10 ONU"
$ws.Range("H11").Value = "3.151"
$ws.Range("I11").Value = "1.909"
$ws.Range("J11").Value = "0"
$ws.Range("K11").Value = "0"
$ws.Range("L11").Value = "0"
$ws.Range("M11").Value = "0"
$ws.Range("N11").Value = "3.151"
$ws.Range("O11").Value = "1.909"
$ws.Range("P11").Value = "0"

# Row 12
$ws.Range("A12").Value = "TQ-17"
$ws.Range("B12").Value = "MERCADO INTERNO"
$ws.Range("C12").Value = "VIBRA PAULINIA"
$ws.Range("D12").Value = "0177/23C"
$ws.Range("E12").Value = "TRANSF  VIBRA CUB P/VIBRA PAULINIA CONF NF 1283140-0 DE09.10.23"
$ws.Range("F12").Value = "OLEO DIESEL S10"
$ws.Range("G12").Value = "OLEO DIESEL S10 ONU"
$ws.Range("H12").Value = "10.136"
$ws.Range("I12").Value = "10.110"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "0"
$ws.Range("L12").Value = "0"
$ws.Range("M12").Value = "0"
$ws.Range("N12").Value = "10.136"
$ws.Range("O12").Value = "10.110"
$ws.Range("P12").Value = "0"

# Row 13
$ws.Range("A13").Value = "TQ-17"
$ws.Range("B13").Value = "MERCADO INTERNO"
$ws.Range("C13").Value = "VIBRA BARUERI"
$ws.Range("D13").Value = "0177/23E"
$ws.Range("E13").Value = "TRANSF  VIBRA CUB X VIBRA BARUERI CONF NF 1283141 09.10.23"
$ws.Range("F13").Value = "OLEO DIESEL S10"
$ws.Range("G13").Value = "OLEO DIESEL S10 ONU"
$ws.Range("H13").Value = "30.435"
$ws.Range("I13").Value = "25.090"
$ws.Range("J13").Value = "0"
$ws.Range("K13").Value = "0"
$ws.Range("L13").Value = "0"
$ws.Range("M13").Value = "0"
$ws.Range("N13").Value = "30.435"
$ws.Range("O13").Value = "25.090"
$ws.Range("P13").Value = "0"

# Row 14
$ws.Range("A14").Value = "TQ-23"
$ws.Range("B14").Value = "MERCADO INTERNO"
$ws.Range("C14").Value = "VIBRA SJ CAMPOS"
$ws.Range("D14").Value = "0053/22E"
$ws.Range("E14").Value = "VIBRA - RIDGEBURY CINDY A - BL RBYK23546293 - 0,8287"
$ws.Range("F14").Value = "OLEO DIESEL S500"
$ws.Range("G14").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H14").Value = "604"
$ws.Range("I14").Value = "501"
$ws.Range("J14").Value = "0"
$ws.Range("K14").Value = "0"
$ws.Range("L14").Value = "0"
$ws.Range("M14").Value = "0"
$ws.Range("N14").Value = "604"
$ws.Range("O14").Value = "501"
$ws.Range("P14").Value = "0"

# Row 15
$ws.Range("A15").Value = "TQ-23"
$ws.Range("B15").Value = "MERCADO INTERNO"
$ws.Range("C15").Value = "VIBRA SAO PAULO"
$ws.Range("D15").Value = "0113/22I"
$ws.Range("E15").Value = "VIBRA - ZANDOLIE - BL 6 - 0,8380"
$ws.Range("F15").Value = "OLEO DIESEL S500"
$ws.Range("G15").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H15").Value = "235"
$ws.Range("I15").Value = "0"
$ws.Range("J15").Value = "0"
$ws.Range("K15").Value = "0"
$ws.Range("L15").Value = "0"
$ws.Range("M15").Value = "0"
$ws.Range("N15").Value = "235"
$ws.Range("O15").Value = "0"
$ws.Range("P15").Value = "0"

# Row 16
$ws.Range("A16").Value = "TQ-23"
$ws.Range("B16").Value = "MERCADO INTERNO"
$ws.Range("C16").Value = "VIBRA RIBEIRAO PRETO"
$ws.Range("D16").Value = "0192/23D"
$ws.Range("E16").Value = "TRANSF V.CUB X R.PRETO CONF NF 12470085-0 DENS:0,8298"
$ws.Range("F16").Value = "OLEO DIESEL S500"
$ws.Range("G16").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H16").Value = "22.000"
$ws.Range("I16").Value = "18.408"
$ws.Range("J16").Value = "0"
$ws.Range("K16").Value = "0"
$ws.Range("L16").Value = "0"
$ws.Range("M16").Value = "0"
$ws.Range("N16").Value = "22.000"
$ws.Range("O16").Value = "18.408"
$ws.Range("P16").Value = "0"

# Row 17
$ws.Range("A17").Value = "TQ-25"
$ws.Range("B17").Value = "MERCADO INTERNO"
$ws.Range("C17").Value = "VIBRA CUBATAO"
$ws.Range("D17").Value = "0228/23"
$ws.Range("E17").Value = "NIMOFAST - ROMEOS - BL: 06 DENS: 0,8281"
$ws.Range("F17").Value = "OLEO DIESEL S500"
$ws.Range("G17").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H17").Value = "1.155.112"
$ws.Range("I17").Value = "955.258"
$ws.Range("J17").Value = "0"
$ws.Range("K17").Value = "0"
$ws.Range("L17").Value = "0"
$ws.Range("M17").Value = "0"
$ws.Range("N17").Value = "1.155.112"
$ws.Range("O17").Value = "955.258"
$ws.Range("P17").Value = "0"

# Row 18
$ws.Range("A18").Value = "TQ-29"
$ws.Range("B18").Value = "MERCADO INTERNO"
$ws.Range("C18").Value = "VIBRA CUBATAO"
$ws.Range("D18").Value = "0228/23"
$ws.Range("E18").Value = "NIMOFAST - ROMEOS - BL: 06 DENS: 0,8281"
$ws.Range("F18").Value = "OLEO DIESEL S500"
$ws.Range("G18").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H18").Value = "578.249"
$ws.Range("I18").Value = "478.202"
$ws.Range("J18").Value = "0"
$ws.Range("K18").Value = "0"
$ws.Range("L18").Value = "0"
$ws.Range("M18").Value = "0"
$ws.Range("N18").Value = "578.249"
$ws.Range("O18").Value = "478.202"
$ws.Range("P18").Value = "0"

# Row 19
$ws.Range("A19").Value = "TQ-38"
$ws.Range("B19").Value = "IMPORTACAO COMUM"
$ws.Range("C19").Value = "VIBRA CUBATAO"
$ws.Range("D19").Value = "0192/23A"
$ws.Range("E19").Value = "TORM ATLANTIC - VIBRA CUBATÃO BL:05 - DENS:0,8305"
$ws.Range("F19").Value = "OLEO DIESEL S500"
$ws.Range("G19").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H19").Value = "1.805"
$ws.Range("I19").Value = "2.000"
$ws.Range("J19").Value = "0"
$ws.Range("K19").Value = "0"
$ws.Range("L19").Value = "0"
$ws.Range("M19").Value = "0"
$ws.Range("N19").Value = "1.805"
$ws.Range("O19").Value = "2.000"
$ws.Range("P19").Value = "0"

# Row 20
$ws.Range("A20").Value = "TQ-38"
$ws.Range("B20").Value = "MERCADO INTERNO"
$ws.Range("C20").Value = "VIBRA SAO PAULO"
$ws.Range("D20").Value = "0192/23C"
$ws.Range("E20").Value = "TRANSF V.CUB X SP CONF NF 1285073-0 DENS:0,8298"
$ws.Range("F20").Value = "OLEO DIESEL S500"
$ws.Range("G20").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H20").Value = "22.696"
$ws.Range("I20").Value = "18.460"
$ws.Range("J20").Value = "0"
$ws.Range("K20").Value = "0"
$ws.Range("L20").Value = "0"
$ws.Range("M20").Value = "0"
$ws.Range("N20").Value = "22.696"
$ws.Range("O20").Value = "18.460"
$ws.Range("P20").Value = "0"

# Row 21
$ws.Range("A21").Value = "TQ-38"
$ws.Range("B21").Value = "MERCADO INTERNO"
$ws.Range("C21").Value = "VIBRA RIBEIRAO PRETO"
$ws.Range("D21").Value = "0192/23D"
$ws.Range("E21").Value = "TRANSF V.CUB X R.PRETO CONF NF 12470085-0 DENS:0,8298"
$ws.Range("F21").Value = "OLEO DIESEL S500"
$ws.Range("G21").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H21").Value = "467.996"
$ws.Range("I21").Value = "388.098"
$ws.Range("J21").Value = "0"
$ws.Range("K21").Value = "0"
$ws.Range("L21").Value = "0"
$ws.Range("M21").Value = "0"
$ws.Range("N21").Value = "467.996"
$ws.Range("O21").Value = "388.098"
$ws.Range("P21").Value = "0"

# Row 22
$ws.Range("A22").Value = "TQ-38"
$ws.Range("B22").Value = "MERCADO INTERNO"
$ws.Range("C22").Value = "VIBRA CUBATAO"
$ws.Range("D22").Value = "0192/23E"
$ws.Range("E22").Value = "TRANSF PAULINIA-CUBATAO CONF NF 3749230-1 DENS:0,8298"
$ws.Range("F22").Value = "OLEO DIESEL S500"
$ws.Range("G22").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H22").Value = "10.112"
$ws.Range("I22").Value = "9.685"
$ws.Range("J22").Value = "0"
$ws.Range("K22").Value = "0"
$ws.Range("L22").Value = "0"
$ws.Range("M22").Value = "0"
$ws.Range("N22").Value = "10.112"
$ws.Range("O22").Value = "9.685"
$ws.Range("P22").Value = "0"

# Row 23
$ws.Range("A23").Value = "TQ-43"
$ws.Range("B23").Value = "MERCADO INTERNO"
$ws.Range("C23").Value = "VIBRA CUBATAO"
$ws.Range("D23").Value = "0228/23"
$ws.Range("E23").Value = "NIMOFAST - ROMEOS - BL: 06 DENS: 0,8281"
$ws.Range("F23").Value = "OLEO DIESEL S500"
$ws.Range("G23").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H23").Value = "2.008.248"
$ws.Range("I23").Value = "1.660.788"
$ws.Range("J23").Value = "0"
$ws.Range("K23").Value = "0"
$ws.Range("L23").Value = "0"
$ws.Range("M23").Value = "0"
$ws.Range("N23").Value = "2.008.248"
$ws.Range("O23").Value = "1.660.788"
$ws.Range("P23").Value = "0"

# Row 24
$ws.Range("A24").Value = "TQ-59"
$ws.Range("B24").Value = "MERCADO INTERNO"
$ws.Range("C24").Value = "VIBRA CUBATAO"
$ws.Range("D24").Value = "0228/23"
$ws.Range("E24").Value = "NIMOFAST - ROMEOS - BL: 06 DENS: 0,8281"
$ws.Range("F24").Value = "OLEO DIESEL S500"
$ws.Range("G24").Value = "ONU 1202, OLEO DIESEL S500"
$ws.Range("H24").Value = "1.253.530"
$ws.Range("I24").Value = "1.036.649"
$ws.Range("J24").Value = "0"
$ws.Range("K24").Value = "0"
$ws.Range("L24").Value = "0"
$ws.Range("M24").Value = "0"
$ws.Range("N24").Value = "1.253.530"
$ws.Range("O24").Value = "1.036.649"
$ws.Range("P24").Value = "0"

# Row 25
$ws.Range("A25").Value = "Z"
$ws.Range("D25").Value = "Total"
$ws.Range("H25").Value = "5.646.865"
$ws.Range("I25").Value = "4.680.239"
$ws.Range("J25").Value = "0"
$ws.Range("K25").Value = "0"
$ws.Range("L25").Value = "0"
$ws.Range("M25").Value = "0"
$ws.Range("N25").Value = "5.646.865"
$ws.Range("O25").Value = "4.680.239"
$ws.Range("P25").Value = "0"
